$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.83171066666667
$ws.Range("H2").Value = 35.495132
$ws.Range("I2").Value = 0.0821640352811125
$ws.Range("J2").Value = 0.08216403528111249
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 18.382477
$ws.Range("N2").Value = 55.147431
$ws.Range("O2").Value = 0.06380158579420245
$ws.Range("P2").Value = 0.06380158579420243
$ws.Range("Q2").Value = 217.4961492006546
$ws.Range("R2").Value = 1957.465342805892
$ws.Range("S2").Value = 0.005242195746185776
$ws.Range("T2").Value = 0.005242195746185774

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.83171066666667
$ws.Range("H3").Value = 35.495132
$ws.Range("I3").Value = 0.0821640352811125
$ws.Range("J3").Value = 0.08216403528111249
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 82.64333833333332
$ws.Range("N3").Value = 247.930015
$ws.Range("O3").Value = 0.2868370808239535
$ws.Range("P3").Value = 0.2868370808239535
$ws.Range("Q3").Value = 977.812067687442
$ws.Range("R3").Value = 8800.308609186979
$ws.Range("S3").Value = 0.02356769202875063
$ws.Range("T3").Value = 0.02356769202875063

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.83171066666667
$ws.Range("H4").Value = 35.495132
$ws.Range("I4").Value = 0.0821640352811125
$ws.Range("J4").Value = 0.08216403528111249
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 91.25099666666667
$ws.Range("N4").Value = 273.75299
$ws.Range("O4").Value = 0.3167123936907314
$ws.Range("P4").Value = 0.3167123936907314
$ws.Range("Q4").Value = 1079.655390604964
$ws.Range("R4").Value = 9716.89851544468
$ws.Range("S4").Value = 0.02602236828917085
$ws.Range("T4").Value = 0.02602236828917084

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.83171066666667
$ws.Range("H5").Value = 35.495132
$ws.Range("I5").Value = 0.0821640352811125
$ws.Range("J5").Value = 0.08216403528111249
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 81.28845566666666
$ws.Range("N5").Value = 243.865367
$ws.Range("O5").Value = 0.2821345773094157
$ws.Range("P5").Value = 0.2821345773094157
$ws.Range("Q5").Value = 961.7814879881604
$ws.Range("R5").Value = 8656.033391893443
$ws.Range("S5").Value = 0.02318131536407259
$ws.Range("T5").Value = 0.02318131536407259

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.83171066666667
$ws.Range("H6").Value = 35.495132
$ws.Range("I6").Value = 0.0821640352811125
$ws.Range("J6").Value = 0.08216403528111249
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.55416966666667
$ws.Range("N6").Value = 43.662509
$ws.Range("O6").Value = 0.0505143623816971
$ws.Range("P6").Value = 0.0505143623816971
$ws.Range("Q6").Value = 172.2007244895764
$ws.Range("R6").Value = 1549.806520406188
$ws.Range("S6").Value = 0.004150463852932663
$ws.Range("T6").Value = 0.004150463852932662

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 42.86866899999999
$ws.Range("H7").Value = 128.606007
$ws.Range("I7").Value = 0.2976968361890019
$ws.Range("J7").Value = 0.2976968361890019
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.382477
$ws.Range("N7").Value = 55.147431
$ws.Range("O7").Value = 0.06380158579420245
$ws.Range("P7").Value = 0.06380158579420243
$ws.Range("Q7").Value = 788.0323219131127
$ws.Range("R7").Value = 7092.290897218016
$ws.Range("S7").Value = 0.01899353023477524
$ws.Range("T7").Value = 0.01899353023477523

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 42.86866899999999
$ws.Range("H8").Value = 128.606007
$ws.Range("I8").Value = 0.2976968361890019
$ws.Range("J8").Value = 0.2976968361890019
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 82.64333833333332
$ws.Range("N8").Value = 247.930015
$ws.Range("O8").Value = 0.2868370808239535
$ws.Range("P8").Value = 0.2868370808239535
$ws.Range("Q8").Value = 3542.809916066677
$ws.Range("R8").Value = 31885.2892446001
$ws.Range("S8").Value = 0.08539049146297997
$ws.Range("T8").Value = 0.08539049146297997

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42.86866899999999
$ws.Range("H9").Value = 128.606007
$ws.Range("I9").Value = 0.2976968361890019
$ws.Range("J9").Value = 0.2976968361890019
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 91.25099666666667
$ws.Range("N9").Value = 273.75299
$ws.Range("O9").Value = 0.3167123936907314
$ws.Range("P9").Value = 0.3167123936907314
$ws.Range("Q9").Value = 3911.808772023436
$ws.Range("R9").Value = 35206.27894821092
$ws.Range("S9").Value = 0.09428427758357634
$ws.Range("T9").Value = 0.09428427758357634

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 42.86866899999999
$ws.Range("H10").Value = 128.606007
$ws.Range("I10").Value = 0.2976968361890019
$ws.Range("J10").Value = 0.2976968361890019
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 81.28845566666666
$ws.Range("N10").Value = 243.865367
$ws.Range("O10").Value = 0.2821345773094157
$ws.Range("P10").Value = 0.2821345773094157
$ws.Range("Q10").Value = 3484.727899495507
$ws.Range("R10").Value = 31362.55109545956
$ws.Range("S10").Value = 0.0839905710445344
$ws.Range("T10").Value = 0.0839905710445344

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 42.86866899999999
$ws.Range("H11").Value = 128.606007
$ws.Range("I11").Value = 0.2976968361890019
$ws.Range("J11").Value = 0.2976968361890019
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 14.55416966666667
$ws.Range("N11").Value = 43.662509
$ws.Range("O11").Value = 0.0505143623816971
$ws.Range("P11").Value = 0.0505143623816971
$ws.Range("Q11").Value = 623.9178820101736
$ws.Range("R11").Value = 5615.260938091562
$ws.Range("S11").Value = 0.01503796586313596
$ws.Range("T11").Value = 0.01503796586313596

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 42.11645766666667
$ws.Range("H12").Value = 126.349373
$ws.Range("I12").Value = 0.2924731859264094
$ws.Range("J12").Value = 0.2924731859264094
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.382477
$ws.Range("N12").Value = 55.147431
$ws.Range("O12").Value = 0.06380158579420245
$ws.Range("P12").Value = 0.06380158579420243
$ws.Range("Q12").Value = 774.2048143789737
$ws.Range("R12").Value = 6967.843329410764
$ws.Range("S12").Value = 0.01866025306438754
$ws.Range("T12").Value = 0.01866025306438753

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 42.11645766666667
$ws.Range("H13").Value = 126.349373
$ws.Range("I13").Value = 0.2924731859264094
$ws.Range("J13").Value = 0.2924731859264094
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 82.64333833333332
$ws.Range("N13").Value = 247.930015
$ws.Range("O13").Value = 0.2868370808239535
$ws.Range("P13").Value = 0.2868370808239535
$ws.Range("Q13").Value = 3480.644660347843
$ws.Range("R13").Value = 31325.80194313059
$ws.Range("S13").Value = 0.08389215487041267
$ws.Range("T13").Value = 0.08389215487041266

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 42.11645766666667
$ws.Range("H14").Value = 126.349373
$ws.Range("I14").Value = 0.2924731859264094
$ws.Range("J14").Value = 0.2924731859264094
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 91.25099666666667
$ws.Range("N14").Value = 273.75299
$ws.Range("O14").Value = 0.3167123936907314
$ws.Range("P14").Value = 0.3167123936907314
$ws.Range("Q14").Value = 3843.168738152808
$ws.Range("R14").Value = 34588.51864337528
$ws.Range("S14").Value = 0.09262988280510746
$ws.Range("T14").Value = 0.09262988280510745

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 42.11645766666667
$ws.Range("H15").Value = 126.349373
$ws.Range("I15").Value = 0.2924731859264094
$ws.Range("J15").Value = 0.2924731859264094
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 81.28845566666666
$ws.Range("N15").Value = 243.865367
$ws.Range("O15").Value = 0.2821345773094157
$ws.Range("P15").Value = 0.2821345773094157
$ws.Range("Q15").Value = 3423.581801873877
$ws.Range("R15").Value = 30812.23621686489
$ws.Range("S15").Value = 0.08251679868568566
$ws.Range("T15").Value = 0.08251679868568564

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 42.11645766666667
$ws.Range("H16").Value = 126.349373
$ws.Range("I16").Value = 0.2924731859264094
$ws.Range("J16").Value = 0.2924731859264094
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 14.55416966666667
$ws.Range("N16").Value = 43.662509
$ws.Range("O16").Value = 0.0505143623816971
$ws.Range("P16").Value = 0.0505143623816971
$ws.Range("Q16").Value = 612.9700706396508
$ws.Range("R16").Value = 5516.730635756858
$ws.Range("S16").Value = 0.01477409650081612
$ws.Range("T16").Value = 0.01477409650081611

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 40.18211366666667
$ws.Range("H17").Value = 120.546341
$ws.Range("I17").Value = 0.2790403431922163
$ws.Range("J17").Value = 0.2790403431922163
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.382477
$ws.Range("N17").Value = 55.147431
$ws.Range("O17").Value = 0.06380158579420245
$ws.Range("P17").Value = 0.06380158579420243
$ws.Range("Q17").Value = 738.6467802888856
$ws.Range("R17").Value = 6647.821022599971
$ws.Range("S17").Value = 0.01780321639622188
$ws.Range("T17").Value = 0.01780321639622188

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 40.18211366666667
$ws.Range("H18").Value = 120.546341
$ws.Range("I18").Value = 0.2790403431922163
$ws.Range("J18").Value = 0.2790403431922163
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 82.64333833333332
$ws.Range("N18").Value = 247.930015
$ws.Range("O18").Value = 0.2868370808239535
$ws.Range("P18").Value = 0.2868370808239535
$ws.Range("Q18").Value = 3320.78401470279
$ws.Range("R18").Value = 29887.05613232511
$ws.Range("S18").Value = 0.08003911747336946
$ws.Range("T18").Value = 0.08003911747336946

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 40.18211366666667
$ws.Range("H19").Value = 120.546341
$ws.Range("I19").Value = 0.2790403431922163
$ws.Range("J19").Value = 0.2790403431922163
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 91.25099666666667
$ws.Range("N19").Value = 273.75299
$ws.Range("O19").Value = 0.3167123936907314
$ws.Range("P19").Value = 0.3167123936907314
$ws.Range("Q19").Value = 3666.657920256621
$ws.Range("R19").Value = 32999.92128230959
$ws.Range("S19").Value = 0.08837553502869
$ws.Range("T19").Value = 0.08837553502869

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 40.18211366666667
$ws.Range("H20").Value = 120.546341
$ws.Range("I20").Value = 0.2790403431922163
$ws.Range("J20").Value = 0.2790403431922163
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 81.28845566666666
$ws.Range("N20").Value = 243.865367
$ws.Range("O20").Value = 0.2821345773094157
$ws.Range("P20").Value = 0.2821345773094157
$ws.Range("Q20").Value = 3266.341965385794
$ws.Range("R20").Value = 29397.07768847215
$ws.Range("S20").Value = 0.07872692927881021
$ws.Range("T20").Value = 0.07872692927881021

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 40.18211366666667
$ws.Range("H21").Value = 120.546341
$ws.Range("I21").Value = 0.2790403431922163
$ws.Range("J21").Value = 0.2790403431922163
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 14.55416966666667
$ws.Range("N21").Value = 43.662509
$ws.Range("O21").Value = 0.0505143623816971
$ws.Range("P21").Value = 0.0505143623816971
$ws.Range("Q21").Value = 584.8172998699521
$ws.Range("R21").Value = 5263.355698829569
$ws.Range("S21").Value = 0.01409554501512474
$ws.Range("T21").Value = 0.01409554501512474

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 7.002139333333333
$ws.Range("H22").Value = 21.006418
$ws.Range("I22").Value = 0.04862559941126002
$ws.Range("J22").Value = 0.04862559941126002
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 18.382477
$ws.Range("N22").Value = 55.147431
$ws.Range("O22").Value = 0.06380158579420245
$ws.Range("P22").Value = 0.06380158579420243
$ws.Range("Q22").Value = 128.7166652457953
$ws.Range("R22").Value = 1158.449987212158
$ws.Range("S22").Value = 0.003102390352632026
$ws.Range("T22").Value = 0.003102390352632026

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 7.002139333333333
$ws.Range("H23").Value = 21.006418
$ws.Range("I23").Value = 0.04862559941126002
$ws.Range("J23").Value = 0.04862559941126002
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 82.64333833333332
$ws.Range("N23").Value = 247.930015
$ws.Range("O23").Value = 0.2868370808239535
$ws.Range("P23").Value = 0.2868370808239535
$ws.Range("Q23").Value = 578.6801699818077
$ws.Range("R23").Value = 5208.121529836269
$ws.Range("S23").Value = 0.01394762498844077
$ws.Range("T23").Value = 0.01394762498844077

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 7.002139333333333
$ws.Range("H24").Value = 21.006418
$ws.Range("I24").Value = 0.04862559941126002
$ws.Range("J24").Value = 0.04862559941126002
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 91.25099666666667
$ws.Range("N24").Value = 273.75299
$ws.Range("O24").Value = 0.3167123936907314
$ws.Range("P24").Value = 0.3167123936907314
$ws.Range("Q24").Value = 638.9521929655356
$ws.Range("R24").Value = 5750.56973668982
$ws.Range("S24").Value = 0.01540032998418678
$ws.Range("T24").Value = 0.01540032998418678

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 7.002139333333333
$ws.Range("H25").Value = 21.006418
$ws.Range("I25").Value = 0.04862559941126002
$ws.Range("J25").Value = 0.04862559941126002
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 81.28845566666666
$ws.Range("N25").Value = 243.865367
$ws.Range("O25").Value = 0.2821345773094157
$ws.Range("P25").Value = 0.2821345773094157
$ws.Range("Q25").Value = 569.1930927694896
$ws.Range("R25").Value = 5122.737834925406
$ws.Range("S25").Value = 0.01371896293631282
$ws.Range("T25").Value = 0.01371896293631282

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 7.002139333333333
$ws.Range("H26").Value = 21.006418
$ws.Range("I26").Value = 0.04862559941126002
$ws.Range("J26").Value = 0.04862559941126002
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 14.55416966666667
$ws.Range("N26").Value = 43.662509
$ws.Range("O26").Value = 0.0505143623816971
$ws.Range("P26").Value = 0.0505143623816971
$ws.Range("Q26").Value = 101.9103238869736
$ws.Range("R26").Value = 917.192914982762
$ws.Range("S26").Value = 0.002456291149687626
$ws.Range("T26").Value = 0.002456291149687626
